$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.921.37'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').Value = '2.671.39'
$ws.Range('E3').Value = '  +2.55%  '
$ws.Range('E4').Value = '  -0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '579.96'
$c.ClearFormats()
$ws.Range('E5').Value = '  +1.14%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '145.11'
$c.ClearFormats()
$ws.Range('E6').Value = '  +1.61%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('E8').Value = '  -0.40%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '6.57'
$c.ClearFormats()
$ws.Range('E9').Value = '  +1.30%  '
$ws.Range('E10').Value = '  +1.60%  '
$ws.Range('E11').Value = '  +4.44%  '
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').Value = '3.143.35'
$ws.Range('E13').Value = '  +2.42%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '25.72'
$c.ClearFormats()
$ws.Range('E14').Value = '  +10.74%  '
$ws.Range('D15').Value = '60.916.39'
$ws.Range('E15').Value = '  +0.70%  '
$ws.Range('E16').Value = '  +1.72%  '
$ws.Range('D17').Value = '2.673.35'
$ws.Range('E18').Value = '  +2.36%  '
$ws.Range('E19').Value = '  +1.61%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '350.88'
$c.ClearFormats()
$ws.Range('E20').Value = '  +1.04%  '
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('E22').Value = '  +0.18%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '0.534'
$c.ClearFormats()
$ws.Range('E23').Value = '  +1.54%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '64.03'
$c.ClearFormats()
$ws.Range('E24').Value = '  +1.26%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('E26').Value = '  +1.76%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '8.12'
$c.ClearFormats()
$ws.Range('E27').Value = '  +4.84%  '
$ws.Range('E28').Value = '  +7.13%  '
$ws.Range('D29').Value = '0.0₃0813'
$ws.Range('E29').Value = '  +2.94%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '6.79'
$c.ClearFormats()
$ws.Range('E30').Value = '  +6.80%  '
$ws.Range('E31').Value = '  +0.13%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '166.96'
$c.ClearFormats()
$ws.Range('E32').Value = '  +3.38%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '19.88'
$c.ClearFormats()
$ws.Range('E33').Value = '  +1.89%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.07'
$c.ClearFormats()
$ws.Range('E34').Value = '  +9.48%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '4.46'
$c.ClearFormats()
$ws.Range('E35').Value = '  +5.50%  '
$ws.Range('E36').Value = '  +7.80%  '
$ws.Range('E37').Value = '  +3.14%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '331.92'
$c.ClearFormats()
$ws.Range('E38').Value = '  +12.48%  '
$ws.Range('E39').Value = '  +4.32%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.881'
$c.ClearFormats()
$ws.Range('E41').Value = '  +4.70%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '20.51'
$c.ClearFormats()
$ws.Range('E42').Value = '  +3.67%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '5.19'
$c.ClearFormats()
$ws.Range('E43').Value = '  +5.17%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '134.65'
$c.ClearFormats()
$ws.Range('E44').Value = '  -1.93%  '
$ws.Range('E45').Value = '  +1.59%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.615'
$c.ClearFormats()
$ws.Range('E47').Value = '  +0.92%  '
$ws.Range('E48').Value = '  +3.19%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E49').Value = '  +0.38%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '20.38'
$c.ClearFormats()
$ws.Range('E50').Value = '  +3.16%  '
$ws.Range('D51').Value = '2.114.39'
$ws.Range('E51').Value = '  +4.34%  '
